# Generate Report for Handback
# Update the timestamp cells that record when the handoff/handback
# report was (re-)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-09-05 11:36:25"

# zh-cn sheet: Correspond Handoff / Handback datetime for the first file
$wsZhCn.Range("H2").Value = "2016-09-05 11:36:20"
$wsZhCn.Range("K2").Value = "2016-09-05 11:36:39"

# de-de sheet: Correspond Handoff / Handback datetime for the first file
$wsDeDe.Range("H2").Value = "2016-09-05 11:36:25"
$wsDeDe.Range("K2").Value = "2016-09-05 11:36:47"
